$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 11833.333
$ws.Range("I32").Value = 9250
$ws.Range("J32").Value = 13125
$ws.Range("K32").Value = 9250
$ws.Range("L32").Value = 13125
$ws.Range("M32").Value = -8924
$ws.Range("N32").Value = -13777
$ws.Range("H62").Value = 33629.242
$ws.Range("I62").Value = 47148.086
$ws.Range("K62").Value = 47148.086
$ws.Range("M62").Value = -46524.086
$ws.Range("H65").Value = 33629.242
$ws.Range("I65").Value = 47148.086
$ws.Range("K65").Value = 235740.43
$ws.Range("M65").Value = -232620.43
$ws.Range("H92").Value = 4303.654
$ws.Range("I92").Value = 5476.95
$ws.Range("J92").Value = 392.66666
$ws.Range("K92").Value = 5476.95
$ws.Range("L92").Value = 392.66666
$ws.Range("M92").Value = -4228.95
$ws.Range("N92").Value = -2888.66666
$ws.Range("H93").Value = 73166.336
$ws.Range("J93").Value = 73166.336
$ws.Range("L93").Value = 73166.336
$ws.Range("N93").Value = -78158.336
$ws.Range("H98").Value = 820.2857
$ws.Range("I98").Value = 815.5333000000001
$ws.Range("J98").Value = 832.1667
$ws.Range("K98").Value = 815.5333000000001
$ws.Range("L98").Value = 832.1667
$ws.Range("M98").Value = 682.4666999999999
$ws.Range("N98").Value = -3828.1667
$ws.Range("H105").Value = 40999.4
$ws.Range("J105").Value = 40999.4
$ws.Range("L105").Value = 40999.4
$ws.Range("N105").Value = -47987.4
$ws.Range("H106").Value = 1810
$ws.Range("I106").Value = 1113.3334
$ws.Range("K106").Value = 1113.3334
$ws.Range("M106").Value = -482.3334
$ws.Range("H116").Value = 8974.546
$ws.Range("I116").Value = 10252.875
$ws.Range("J116").Value = 5565.6665
$ws.Range("K116").Value = 10252.875
$ws.Range("L116").Value = 5565.6665
$ws.Range("M116").Value = -6810.875
$ws.Range("N116").Value = -12449.6665
$ws.Range("H122").Value = 820.2857
$ws.Range("I122").Value = 815.5333000000001
$ws.Range("J122").Value = 832.1667
$ws.Range("K122").Value = 2446.5999
$ws.Range("L122").Value = 2496.5001
$ws.Range("M122").Value = 3.400099999999838
$ws.Range("N122").Value = -7396.5001
$ws.Range("H132").Value = 62201.21
$ws.Range("I132").Value = 40024.617
$ws.Range("K132").Value = 120073.851
$ws.Range("M132").Value = -117543.851
$ws.Range("H141").Value = 1498.3334
$ws.Range("I141").Value = 1297.5
$ws.Range("K141").Value = 3892.5
$ws.Range("M141").Value = 1287.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 398.4
$ws.Range("I5").Value = 255.88889
$ws.Range("K5").Value = 255.88889
$ws.Range("M5").Value = -143.88889
$ws.Range("H97").Value = 3294
$ws.Range("I97").Value = 3343
$ws.Range("K97").Value = 3343
$ws.Range("M97").Value = -2847
$ws.Range("H102").Value = 4654.8887
$ws.Range("I102").Value = 4482.5
$ws.Range("J102").Value = 4999.6665
$ws.Range("K102").Value = 4482.5
$ws.Range("L102").Value = 4999.6665
$ws.Range("M102").Value = -2860.5
$ws.Range("N102").Value = -8243.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 398.4
$ws.Range("I4").Value = 255.88889
$ws.Range("K4").Value = 255.88889
$ws.Range("M4").Value = -140.88889
$ws.Range("H80").Value = 1619.375
$ws.Range("J80").Value = 1827.1666
$ws.Range("L80").Value = 1827.1666
$ws.Range("N80").Value = -3823.1666
$ws.Range("H83").Value = 1619.375
$ws.Range("J83").Value = 1827.1666
$ws.Range("L83").Value = 9135.833000000001
$ws.Range("N83").Value = -19119.833
$ws.Range("H88").Value = 20028000
$ws.Range("I88").Value = 20000
$ws.Range("J88").Value = 25030000
$ws.Range("K88").Value = 20000
$ws.Range("L88").Value = 25030000
$ws.Range("M88").Value = -19594
$ws.Range("N88").Value = -25030812
$ws.Range("H91").Value = 20028000
$ws.Range("I91").Value = 20000
$ws.Range("J91").Value = 25030000
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 25030000
$ws.Range("M91").Value = -18596
$ws.Range("N91").Value = -25032808
$ws.Range("H103").Value = 23868.166
$ws.Range("J103").Value = 23868.166
$ws.Range("L103").Value = 23868.166
$ws.Range("N103").Value = -26212.166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2872.5293
$ws.Range("I99").Value = 2872.5293
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2872.5293
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1374.5293
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 2872.5293
$ws.Range("I126").Value = 2872.5293
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8617.5879
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6147.5879
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 29324.834
$ws.Range("I87").Value = 28649.666
$ws.Range("K87").Value = 85948.99800000001
$ws.Range("M87").Value = -84700.99800000001
$ws.Range("H90").Value = 29324.834
$ws.Range("I90").Value = 28649.666
$ws.Range("K90").Value = 257846.994
$ws.Range("M90").Value = -251606.994
$ws.Range("H99").Value = 9833.333000000001
$ws.Range("I99").Value = 3500
$ws.Range("K99").Value = 10500
$ws.Range("M99").Value = -8254
$ws.Range("H131").Value = 30852.543
$ws.Range("I131").Value = 250498.5
$ws.Range("J131").Value = 2511.1292
$ws.Range("K131").Value = 751495.5
$ws.Range("L131").Value = 7533.3876
$ws.Range("M131").Value = -746455.5
$ws.Range("N131").Value = -17613.3876

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1167
$ws.Range("I97").Value = 986.2857
$ws.Range("J97").Value = 1420
$ws.Range("K97").Value = 986.2857
$ws.Range("L97").Value = 1420
$ws.Range("M97").Value = -490.2857
$ws.Range("N97").Value = -2412
$ws.Range("H132").Value = 388310
$ws.Range("I132").Value = 628685.4
$ws.Range("K132").Value = 1886056.2
$ws.Range("M132").Value = -1883526.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4326.027
$ws.Range("I16").Value = 1322
$ws.Range("K16").Value = 1322
$ws.Range("M16").Value = -1152
$ws.Range("H43").Value = 488433.75
$ws.Range("I43").Value = 1516000
$ws.Range("K43").Value = 1516000
$ws.Range("M43").Value = -1515807
$ws.Range("H55").Value = 264.3
$ws.Range("I55").Value = 281.69232
$ws.Range("K55").Value = 281.69232
$ws.Range("M55").Value = -108.69232
$ws.Range("H93").Value = 2888.238
$ws.Range("I93").Value = 2773.375
$ws.Range("J93").Value = 3255.8
$ws.Range("K93").Value = 2773.375
$ws.Range("L93").Value = 3255.8
$ws.Range("M93").Value = -1525.375
$ws.Range("N93").Value = -5751.8
$ws.Range("H127").Value = 71500
$ws.Range("J127").Value = 71500
$ws.Range("L127").Value = 71500
$ws.Range("N127").Value = -81420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 34999.668
$ws.Range("J92").Value = 37000
$ws.Range("L92").Value = 37000
$ws.Range("N92").Value = -41992
